$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data appended after the last existing row (row 20 -> new row 21)
# Force the date column to be treated as literal text (matching the other
# date cells in column A) instead of being auto-parsed into a date serial.
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "09/22/2025"
$ws.Range("A21").Style = $ws.Range("A2").Style

$ws.Range("B21").Value = 0.1347731634554507
$ws.Range("C21").Value = 0.8652268365445493
